$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

$zhTargetTime = "2016-08-19 08:58:48"
$deTargetTime = "2016-08-19 08:58:55"

$mdUrl1 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/458acc0c33c389eccebc510a39b49dba0e7766f5/e2e/35a68a43-d5c4-4db8-929a-cc1227cc0aad.md"
$mdUrl2 = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/458acc0c33c389eccebc510a39b49dba0e7766f5/e2e/4343d8c3-499b-464b-aef7-5da4f60f5d67.md"
$mdName1 = "35a68a43-d5c4-4db8-929a-cc1227cc0aad.md"
$mdName2 = "4343d8c3-499b-464b-aef7-5da4f60f5d67.md"

$hyperlinkUnderline = $true
$hyperlinkColor = 15570276

# ---------------------------------------------------------------------------
# Overview sheet: status text + widened zh-cn / de-de columns
# ---------------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = $statusHandedBack
$wsOverview.Range("F2").Value = $statusHandedBack
$wsOverview.Range("E3").Value = $statusHandedBack
$wsOverview.Range("F3").Value = $statusHandedBack
$wsOverview.Columns.Item(5).ColumnWidth = 29.9777047293527
$wsOverview.Columns.Item(6).ColumnWidth = 29.9777047293527

# ---------------------------------------------------------------------------
# zh-cn sheet: status text, handback file / datetime, widened columns
# ---------------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("C2").Value = $statusHandedBack
$wsZh.Range("C3").Value = $statusHandedBack

$wsZh.Range("J2").Value = "35a68a43-d5c4-4db8-929a-cc1227cc0aad.a817a7eb539eff2754ffba004f22c1f6b0b9a376.zh-cn.xlf"
$wsZh.Range("K2").Value = $zhTargetTime
$wsZh.Range("J3").Value = "4343d8c3-499b-464b-aef7-5da4f60f5d67.2913515faef2a4266c41c071e21925fc7c847a5c.zh-cn.xlf"
$wsZh.Range("K3").Value = $zhTargetTime

$wsZh.Hyperlinks.Delete()
$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdUrl1, "", "", $mdName1)
$wsZh.Hyperlinks.Add($wsZh.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdUrl2, "", "", $mdName2)
$wsZh.Hyperlinks.Add($wsZh.Range("I3"), $mdUrl2, "", "", $mdName2)

$wsZh.Range("I2").Font.Underline = $hyperlinkUnderline
$wsZh.Range("I2").Font.Color = $hyperlinkColor
$wsZh.Range("I3").Font.Underline = $hyperlinkUnderline
$wsZh.Range("I3").Font.Color = $hyperlinkColor

$wsZh.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsZh.Columns.Item(9).ColumnWidth = 40
$wsZh.Columns.Item(10).ColumnWidth = 40

# ---------------------------------------------------------------------------
# de-de sheet: status text, handback file / datetime, widened columns
# ---------------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("C2").Value = $statusHandedBack
$wsDe.Range("C3").Value = $statusHandedBack

$wsDe.Range("J2").Value = "35a68a43-d5c4-4db8-929a-cc1227cc0aad.a817a7eb539eff2754ffba004f22c1f6b0b9a376.de-de.xlf"
$wsDe.Range("K2").Value = $deTargetTime
$wsDe.Range("J3").Value = "4343d8c3-499b-464b-aef7-5da4f60f5d67.2913515faef2a4266c41c071e21925fc7c847a5c.de-de.xlf"
$wsDe.Range("K3").Value = $deTargetTime

$wsDe.Hyperlinks.Delete()
$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdUrl1, "", "", $mdName1)
$wsDe.Hyperlinks.Add($wsDe.Range("I2"), $mdUrl1, "", "", $mdName1)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdUrl2, "", "", $mdName2)
$wsDe.Hyperlinks.Add($wsDe.Range("I3"), $mdUrl2, "", "", $mdName2)

$wsDe.Range("I2").Font.Underline = $hyperlinkUnderline
$wsDe.Range("I2").Font.Color = $hyperlinkColor
$wsDe.Range("I3").Font.Underline = $hyperlinkUnderline
$wsDe.Range("I3").Font.Color = $hyperlinkColor

$wsDe.Columns.Item(3).ColumnWidth = 29.9777047293527
$wsDe.Columns.Item(9).ColumnWidth = 40
$wsDe.Columns.Item(10).ColumnWidth = 40

Write-Host "Handback report generated."
